$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 113
$ws.Cells.Item(113, 1).NumberFormat = "@"
$ws.Cells.Item(113, 1).Value = "2025-07-18"
$ws.Cells.Item(113, 1).Style = "Normal"
$ws.Cells.Item(113, 2).Value = "Sarpsborg 08 FF"
$ws.Cells.Item(113, 3).Value = "Rosenborg"
$ws.Cells.Item(113, 4).Value = 2
$ws.Cells.Item(113, 5).Value = 2
$ws.Cells.Item(113, 6).Value = 1342294
$ws.Cells.Item(113, 7).Value = 4
$ws.Cells.Item(113, 8).Value = 1
$ws.Cells.Item(113, 9).Value = 2
$ws.Cells.Item(113, 10).Value = 5
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 0
$ws.Cells.Item(113, 14).Value = 0
$ws.Cells.Item(113, 15).Value = 2
$ws.Cells.Item(113, 16).Value = 2
$ws.Cells.Item(113, 17).NumberFormat = "@"
$ws.Cells.Item(113, 17).Value = "52%"
$ws.Cells.Item(113, 17).Style = "Normal"
$ws.Cells.Item(113, 18).NumberFormat = "@"
$ws.Cells.Item(113, 18).Value = "48%"
$ws.Cells.Item(113, 18).Style = "Normal"
$ws.Cells.Item(113, 19).Value = "E"

# Row 114
$ws.Cells.Item(114, 1).NumberFormat = "@"
$ws.Cells.Item(114, 1).Value = "2025-07-19"
$ws.Cells.Item(114, 1).Style = "Normal"
$ws.Cells.Item(114, 2).Value = "KFUM Oslo"
$ws.Cells.Item(114, 3).Value = "Brann"
$ws.Cells.Item(114, 4).Value = 2
$ws.Cells.Item(114, 5).Value = 0
$ws.Cells.Item(114, 6).Value = 1342291
$ws.Cells.Item(114, 7).Value = 4
$ws.Cells.Item(114, 8).Value = 11
$ws.Cells.Item(114, 9).Value = 2
$ws.Cells.Item(114, 10).Value = 5
$ws.Cells.Item(114, 11).Value = 1
$ws.Cells.Item(114, 12).Value = 1
$ws.Cells.Item(114, 13).Value = 0
$ws.Cells.Item(114, 14).Value = 0
$ws.Cells.Item(114, 15).Value = 2
$ws.Cells.Item(114, 16).Value = 0
$ws.Cells.Item(114, 17).NumberFormat = "@"
$ws.Cells.Item(114, 17).Value = "45%"
$ws.Cells.Item(114, 17).Style = "Normal"
$ws.Cells.Item(114, 18).NumberFormat = "@"
$ws.Cells.Item(114, 18).Value = "55%"
$ws.Cells.Item(114, 18).Style = "Normal"
$ws.Cells.Item(114, 19).Value = "L"

# Row 115
$ws.Cells.Item(115, 1).NumberFormat = "@"
$ws.Cells.Item(115, 1).Value = "2025-07-19"
$ws.Cells.Item(115, 1).Style = "Normal"
$ws.Cells.Item(115, 2).Value = "Molde"
$ws.Cells.Item(115, 3).Value = "Stromsgodset"
$ws.Cells.Item(115, 4).Value = 4
$ws.Cells.Item(115, 5).Value = 1
$ws.Cells.Item(115, 6).Value = 1342292
$ws.Cells.Item(115, 7).Value = 6
$ws.Cells.Item(115, 8).Value = 5
$ws.Cells.Item(115, 9).Value = 2
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 11).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 13).Value = 0
$ws.Cells.Item(115, 14).Value = 0
$ws.Cells.Item(115, 15).Value = 4
$ws.Cells.Item(115, 16).Value = 1
$ws.Cells.Item(115, 17).NumberFormat = "@"
$ws.Cells.Item(115, 17).Value = "62%"
$ws.Cells.Item(115, 17).Style = "Normal"
$ws.Cells.Item(115, 18).NumberFormat = "@"
$ws.Cells.Item(115, 18).Value = "38%"
$ws.Cells.Item(115, 18).Style = "Normal"
$ws.Cells.Item(115, 19).Value = "L"

# Row 116
$ws.Cells.Item(116, 1).NumberFormat = "@"
$ws.Cells.Item(116, 1).Value = "2025-07-19"
$ws.Cells.Item(116, 1).Style = "Normal"
$ws.Cells.Item(116, 2).Value = "Viking"
$ws.Cells.Item(116, 3).Value = "Bodo/Glimt"
$ws.Cells.Item(116, 4).Value = 2
$ws.Cells.Item(116, 5).Value = 4
$ws.Cells.Item(116, 6).Value = 1342296
$ws.Cells.Item(116, 7).Value = 4
$ws.Cells.Item(116, 8).Value = 4
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 3
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = 0
$ws.Cells.Item(116, 14).Value = 0
$ws.Cells.Item(116, 15).Value = 2
$ws.Cells.Item(116, 16).Value = 4
$ws.Cells.Item(116, 17).NumberFormat = "@"
$ws.Cells.Item(116, 17).Value = "36%"
$ws.Cells.Item(116, 17).Style = "Normal"
$ws.Cells.Item(116, 18).NumberFormat = "@"
$ws.Cells.Item(116, 18).Value = "64%"
$ws.Cells.Item(116, 18).Style = "Normal"
$ws.Cells.Item(116, 19).Value = "V"

# Row 117
$ws.Cells.Item(117, 1).NumberFormat = "@"
$ws.Cells.Item(117, 1).Value = "2025-07-20"
$ws.Cells.Item(117, 1).Style = "Normal"
$ws.Cells.Item(117, 2).Value = "Valerenga"
$ws.Cells.Item(117, 3).Value = "Haugesund"
$ws.Cells.Item(117, 4).Value = 3
$ws.Cells.Item(117, 5).Value = 0
$ws.Cells.Item(117, 6).Value = 1342297
$ws.Cells.Item(117, 7).Value = 10
$ws.Cells.Item(117, 8).Value = 7
$ws.Cells.Item(117, 9).Value = 1
$ws.Cells.Item(117, 10).Value = 3
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 13).Value = 0
$ws.Cells.Item(117, 14).Value = 0
$ws.Cells.Item(117, 15).Value = 3
$ws.Cells.Item(117, 16).Value = 0
$ws.Cells.Item(117, 17).NumberFormat = "@"
$ws.Cells.Item(117, 17).Value = "49%"
$ws.Cells.Item(117, 17).Style = "Normal"
$ws.Cells.Item(117, 18).NumberFormat = "@"
$ws.Cells.Item(117, 18).Value = "51%"
$ws.Cells.Item(117, 18).Style = "Normal"
$ws.Cells.Item(117, 19).Value = "L"

# Row 118
$ws.Cells.Item(118, 1).NumberFormat = "@"
$ws.Cells.Item(118, 1).Value = "2025-07-20"
$ws.Cells.Item(118, 1).Style = "Normal"
$ws.Cells.Item(118, 2).Value = "Tromso"
$ws.Cells.Item(118, 3).Value = "Bryne"
$ws.Cells.Item(118, 4).Value = 3
$ws.Cells.Item(118, 5).Value = 1
$ws.Cells.Item(118, 6).Value = 1342295
$ws.Cells.Item(118, 7).Value = 5
$ws.Cells.Item(118, 8).Value = 3
$ws.Cells.Item(118, 9).Value = 1
$ws.Cells.Item(118, 10).Value = 3
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 13).Value = 0
$ws.Cells.Item(118, 14).Value = 0
$ws.Cells.Item(118, 15).Value = 3
$ws.Cells.Item(118, 16).Value = 1
$ws.Cells.Item(118, 17).NumberFormat = "@"
$ws.Cells.Item(118, 17).Value = "52%"
$ws.Cells.Item(118, 17).Style = "Normal"
$ws.Cells.Item(118, 18).NumberFormat = "@"
$ws.Cells.Item(118, 18).Value = "48%"
$ws.Cells.Item(118, 18).Style = "Normal"
$ws.Cells.Item(118, 19).Value = "L"

# Row 119
$ws.Cells.Item(119, 1).NumberFormat = "@"
$ws.Cells.Item(119, 1).Value = "2025-07-20"
$ws.Cells.Item(119, 1).Style = "Normal"
$ws.Cells.Item(119, 2).Value = "Sandefjord"
$ws.Cells.Item(119, 3).Value = "Kristiansund BK"
$ws.Cells.Item(119, 4).Value = 6
$ws.Cells.Item(119, 5).Value = 0
$ws.Cells.Item(119, 6).Value = 1342293
$ws.Cells.Item(119, 7).Value = 6
$ws.Cells.Item(119, 8).Value = 7
$ws.Cells.Item(119, 9).Value = 1
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 11).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 13).Value = 0
$ws.Cells.Item(119, 14).Value = 0
$ws.Cells.Item(119, 15).Value = 6
$ws.Cells.Item(119, 16).Value = 0
$ws.Cells.Item(119, 17).NumberFormat = "@"
$ws.Cells.Item(119, 17).Value = "70%"
$ws.Cells.Item(119, 17).Style = "Normal"
$ws.Cells.Item(119, 18).NumberFormat = "@"
$ws.Cells.Item(119, 18).Value = "30%"
$ws.Cells.Item(119, 18).Style = "Normal"
$ws.Cells.Item(119, 19).Value = "L"

# Row 120
$ws.Cells.Item(120, 1).NumberFormat = "@"
$ws.Cells.Item(120, 1).Value = "2025-07-20"
$ws.Cells.Item(120, 1).Style = "Normal"
$ws.Cells.Item(120, 2).Value = "Ham-Kam"
$ws.Cells.Item(120, 3).Value = "Fredrikstad"
$ws.Cells.Item(120, 4).Value = 1
$ws.Cells.Item(120, 5).Value = 1
$ws.Cells.Item(120, 6).Value = 1342290
$ws.Cells.Item(120, 7).Value = 9
$ws.Cells.Item(120, 8).Value = 4
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 1
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 12).Value = 0
$ws.Cells.Item(120, 13).Value = 0
$ws.Cells.Item(120, 14).Value = 0
$ws.Cells.Item(120, 15).Value = 1
$ws.Cells.Item(120, 16).Value = 1
$ws.Cells.Item(120, 17).NumberFormat = "@"
$ws.Cells.Item(120, 17).Value = "43%"
$ws.Cells.Item(120, 17).Style = "Normal"
$ws.Cells.Item(120, 18).NumberFormat = "@"
$ws.Cells.Item(120, 18).Value = "57%"
$ws.Cells.Item(120, 18).Style = "Normal"
$ws.Cells.Item(120, 19).Value = "E"

